$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." spread across several runs
# ("Versi" / "on" / " 2" / "."), with a spell-check proofErr pair around
# "Version" and a _GoBack bookmark between " 2" and ".".
# Target: "Version 1." with "Version" in a single run, the spellEnd
# proofErr right after it, " 1." in a single run, and the (now empty)
# _GoBack bookmark placed after that run (the trailing "." run is gone).

$p1 = $d.Paragraphs(1)

# Range covering the paragraph's content but excluding the trailing
# paragraph-mark character, so InsertXML only touches the run content.
$r = $d.Range($p1.Range.Start, $p1.Range.End - 1)

$xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t>Version</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 1.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></pkg:xmlData>'

$r.InsertXML($xml)
